# Apply Natmi Cd5l-Cd5 update: revise M1 sending-cluster rows and add new
# M2 sending-cluster rows, per Dr Hou advice.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @("M1", "Cd5l", "Cd5", "ECs", 3, 1, 6.795051999999999, 20.385156, 0.6049503343582541, 0.6049503343582541, 3, 1, 0.7920173333333332, 2.376052, 0.1802777704012534, 0.1802777704012534, 5.381798964901332, 48.43619068411199, 0.1090590974815988, 0.1090590974815988)
  3 = @("M1", "Cd5l", "Cd5", "FAPs", 3, 1, 6.795051999999999, 20.385156, 0.6049503343582541, 0.6049503343582541, 3, 1, 1.48869, 4.46607, 0.3388533340414797, 0.3388533340414797, 10.11572596188, 91.04153365692, 0.2049894377268023, 0.2049894377268023)
  4 = @("M1", "Cd5l", "Cd5", "M1", 3, 1, 6.795051999999999, 20.385156, 0.6049503343582541, 0.6049503343582541, 3, 1, 0.4248016666666667, 1.274405, 0.0966927036900747, 0.0966927036900747, 2.886549414686666, 25.97894473218, 0.05849428342731428, 0.05849428342731428)
  5 = @("M1", "Cd5l", "Cd5", "M2", 3, 1, 6.795051999999999, 20.385156, 0.6049503343582541, 0.6049503343582541, 3, 1, 1.426577666666667, 4.279733, 0.3247154200129742, 0.3247154200129743, 9.693669427038666, 87.243024843348, 0.1964367019081296, 0.1964367019081297)
  6 = @("M1", "Cd5l", "Cd5", "sCs", 3, 1, 6.795051999999999, 20.385156, 0.6049503343582541, 0.6049503343582541, 3, 1, 0.26123, 0.7836900000000001, 0.05946077185421796, 0.05946077185421797, 1.77507143396, 15.97564290564, 0.03597081381440902, 0.03597081381440902)
  7 = @("M2", "Cd5l", "Cd5", "ECs", 2, 0.6666666666666666, 4.437361, 13.312083, 0.3950496656417459, 0.3950496656417459, 3, 1, 0.7920173333333332, 2.376052, 0.1802777704012534, 0.1802777704012534, 3.514466826257333, 31.63020143631599, 0.07121867291965459, 0.0712186729196546)
  8 = @("M2", "Cd5l", "Cd5", "FAPs", 2, 0.6666666666666666, 4.437361, 13.312083, 0.3950496656417459, 0.3950496656417459, 3, 1, 1.48869, 4.46607, 0.3388533340414797, 0.3388533340414797, 6.60585494709, 59.45269452381, 0.1338638963146774, 0.1338638963146774)
  9 = @("M2", "Cd5l", "Cd5", "M1", 2, 0.6666666666666666, 4.437361, 13.312083, 0.3950496656417459, 0.3950496656417459, 3, 1, 0.4248016666666667, 1.274405, 0.0966927036900747, 0.0966927036900747, 1.884998348401667, 16.964985135615, 0.03819842026276042, 0.03819842026276042)
  10 = @("M2", "Cd5l", "Cd5", "M2", 2, 0.6666666666666666, 4.437361, 13.312083, 0.3950496656417459, 0.3950496656417459, 3, 1, 1.426577666666667, 4.279733, 0.3247154200129742, 0.3247154200129743, 6.330240101537667, 56.972160913839, 0.1282787181048446, 0.1282787181048446)
  11 = @("M2", "Cd5l", "Cd5", "sCs", 2, 0.6666666666666666, 4.437361, 13.312083, 0.3950496656417459, 0.3950496656417459, 3, 1, 0.26123, 0.7836900000000001, 0.05946077185421796, 0.05946077185421797, 1.15917181403, 10.43254632627, 0.02348995803980894, 0.02348995803980895)
}

foreach ($r in $data.Keys) {
  $vals = $data[$r]
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $col = [char](65 + $i)
    $ws.Range("$col$r").Value = $vals[$i]
  }
}
